# Update cryptocurrency price/volume data on the worksheet.
# Data-driven: each entry is Row -> { Column -> NewValue }
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Cell, $Text) {
    # Force the cell to remain plain text, matching the original inline-string
    # cell type, even when the new value looks like a number (e.g. "1.01").
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

$updates = @{
    2 = @{ "D" = "42.676.14"; "E" = "  -0.28%  " }
    3 = @{ "D" = "2.303.30"; "E" = "  -0.43%  " }
    4 = @{ "D" = "1.01"; "E" = "  +0.24%  " }
    5 = @{ "D" = "310.26"; "E" = "  -2.59%  " }
    6 = @{ "D" = "104.53"; "E" = "  -0.06%  " }
    7 = @{ "E" = "  -0.83%  " }
    8 = @{ "E" = "  +0.14%  " }
    9 = @{ "D" = "0.605"; "E" = "  -0.56%  " }
    10 = @{ "D" = "39.58"; "E" = "  -1.57%  " }
    11 = @{ "D" = "0.0903"; "E" = "  -0.58%  " }
    12 = @{ "D" = "8.27"; "E" = "  -3.38%  " }
    13 = @{ "D" = "0.107"; "E" = "  +0.33%  " }
    14 = @{ "D" = "0.989"; "E" = "  +1.20%  " }
    15 = @{ "D" = "2.783.65"; "E" = "  +4.49%  " }
    16 = @{ "D" = "15.35"; "E" = "  -0.53%  " }
    17 = @{ "D" = "2.298.42"; "E" = "  -0.64%  " }
    18 = @{ "D" = "42.865.95"; "E" = "  +0.27%  " }
    19 = @{ "D" = "7.31"; "E" = "  -2.92%  " }
    20 = @{ "E" = "  -1.29%  " }
    21 = @{ "D" = "13.54"; "E" = "  +1.97%  " }
    22 = @{ "D" = "73.37"; "E" = "  -0.64%  " }
    23 = @{ "E" = "  -2.81%  " }
    24 = @{ "D" = "267.85"; "E" = "  -0.35%  " }
    25 = @{ "D" = "2.20"; "E" = "  -2.67%  " }
    26 = @{ "E" = "  +0.59%  " }
    27 = @{ "D" = "10.88"; "E" = "  -0.19%  " }
    28 = @{ "D" = "7.23"; "E" = "  +15.89%  " }
    29 = @{ "D" = "2.30"; "E" = "  -1.24%  " }
    30 = @{ "D" = "22.30"; "E" = "  -1.85%  " }
    31 = @{ "D" = "36.17"; "E" = "  -4.88%  " }
    32 = @{ "D" = "165.00"; "E" = "  -0.47%  " }
    33 = @{ "D" = "0.0858"; "E" = "  -3.62%  " }
    34 = @{ "E" = "  -1.22%  " }
    35 = @{ "E" = "  +2.27%  " }
    36 = @{ "E" = "  -3.74%  " }
    37 = @{ "D" = "4.53"; "E" = "  -1.87%  " }
    38 = @{ "E" = "  -1.35%  " }
    39 = @{ "D" = "2.81"; "E" = "  +1.86%  " }
    40 = @{ "D" = "3.62"; "E" = "  -2.78%  " }
    41 = @{ "D" = "109.49"; "E" = "  +10.73%  " }
    42 = @{ "D" = "1.58"; "E" = "  -0.73%  " }
    43 = @{ "D" = "70.79"; "E" = "  +0.65%  " }
    44 = @{ "B" = "FirstDigitalUSD"; "C" = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"; "D" = "1.01"; "E" = "  +0.22%  " }
    45 = @{ "B" = "Algorand"; "C" = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; "D" = "0.226"; "E" = "  +0.18%  " }
    46 = @{ "D" = "12.32"; "E" = "  -0.77%  " }
    47 = @{ "D" = "1.723.74"; "E" = "  +6.43%  " }
    48 = @{ "D" = "110.68"; "E" = "  -3.75%  " }
    49 = @{ "D" = "77.51"; "E" = "  -5.84%  " }
    50 = @{ "D" = "5.15"; "E" = "  -3.10%  " }
    51 = @{ "D" = "8.63"; "E" = "  -3.01%  " }
}

foreach ($rowNum in $updates.Keys) {
    $rowData = $updates[$rowNum]
    foreach ($col in $rowData.Keys) {
        $cellRef = "$col$rowNum"
        Set-TextValue $ws.Range($cellRef) $rowData[$col]
    }
}

Write-Host "Updated $($updates.Count) rows."